$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 10:25"

# 2) Rusia (row 7) - numbers refreshed, no rank change
$ws.Range("B7").Value = 789190
$ws.Range("C7").Value = 5862
$ws.Range("D7").Value = 572053
$ws.Range("E7").Value = 204392
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 165
$ws.Range("H7").Value = 12745

# 3) Kuwait / Ucrania swap ranks (row 38 / row 39) with refreshed numbers
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 60995
$ws.Range("C38").Value = 829
$ws.Range("D38").Value = 33172
$ws.Range("E38").Value = 26289
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 16
$ws.Range("H38").Value = 1534

$ws.Range("A39").Value = "Kuwait"
$ws.Range("B39").Value = 60434
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 50919
$ws.Range("E39").Value = 9103
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 412

# 4) Israel (row 43) - numbers refreshed, no rank change
$ws.Range("B43").Value = 54663
$ws.Range("C43").Value = 621
$ws.Range("D43").Value = 22920
$ws.Range("E43").Value = 31313
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 430

# 5) Singapur (row 46) - numbers refreshed, no rank change
$ws.Range("B46").Value = 48744
$ws.Range("C46").Value = 310
$ws.Range("D46").Value = 44584
$ws.Range("E46").Value = 4133
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 27

# 6) Moldavia (row 63) - numbers refreshed, no rank change
$ws.Range("B63").Value = 21442
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 14856
$ws.Range("E63").Value = 5878
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 708

# 7) Estonia / Eslovaquia swap ranks (row 121 / row 122) with refreshed numbers
$ws.Range("A121").Value = "Eslovaquia"
$ws.Range("B121").Value = 2058
$ws.Range("C121").Value = 37
$ws.Range("D121").Value = 1556
$ws.Range("E121").Value = 474
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 28

$ws.Range("A122").Value = "Estonia"
$ws.Range("B122").Value = 2025
$ws.Range("C122").Value = 3
$ws.Range("D122").Value = 1912
$ws.Range("E122").Value = 44
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 69

# 8) Lituania (row 126) - numbers refreshed, no rank change
$ws.Range("B126").Value = 1951
$ws.Range("C126").Value = 2
$ws.Range("D126").Value = 1607
$ws.Range("E126").Value = 264
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 80

# 9) Islas Malvinas / Groenlandia swap ranks (row 210 / row 211) - tied values, no number change
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
